# Auto-generated Excel COM-interop script
# Applies targeted cell value updates/additions/removals across 8 worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 278.30768
$ws.Range("I2").Value = 175.5
$ws.Range("J2").Value = 324
$ws.Range("K2").Value = 175.5
$ws.Range("L2").Value = 324
$ws.Range("M2").Value = -62.5
$ws.Range("N2").Value = -550
$ws.Range("H9").Value = 7816.077
$ws.Range("I9").Value = 11218.777
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 11218.777
$ws.Range("L9").Value = 160
$ws.Range("M9").Value = -11049.777
$ws.Range("N9").Value = -498
$ws.Range("H33").Value = 343.77777
$ws.Range("I33").Value = 184.85715
$ws.Range("J33").Value = 900
$ws.Range("K33").Value = 184.85715
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = 44.14285000000001
$ws.Range("N33").Value = -1358
$ws.Range("H43").Value = 2134.75
$ws.Range("I43").Value = 1920
$ws.Range("J43").Value = 2492.6667
$ws.Range("K43").Value = 1920
$ws.Range("L43").Value = 2492.6667
$ws.Range("M43").Value = -1851
$ws.Range("N43").Value = -2630.6667
$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 35000
$ws.Range("N134").Value = -45140
$ws.Range("H137").Value = 3743.6667
$ws.Range("I137").Value = 3333
$ws.Range("J137").Value = 3949
$ws.Range("K137").Value = 9999
$ws.Range("L137").Value = 11847
$ws.Range("M137").Value = -7449
$ws.Range("N137").Value = -16947
$ws.Range("H141").Value = 11834
$ws.Range("I141").Value = 5225
$ws.Range("K141").Value = 15675
$ws.Range("M141").Value = -10495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6559.8
$ws.Range("I61").Value = 5800
$ws.Range("J61").Value = 7699.5
$ws.Range("K61").Value = 5800
$ws.Range("L61").Value = 7699.5
$ws.Range("M61").Value = -5588
$ws.Range("N61").Value = -8123.5
$ws.Range("H63").Value = 2379.4
$ws.Range("J63").Value = 1998.3334
$ws.Range("L63").Value = 1998.3334
$ws.Range("N63").Value = -3370.3334
$ws.Range("H66").Value = 2379.4
$ws.Range("J66").Value = 1998.3334
$ws.Range("L66").Value = 9991.666999999999
$ws.Range("N66").Value = -16855.667
$ws.Range("H74").Value = 2015.7333
$ws.Range("I74").Value = 1191.6316
$ws.Range("K74").Value = 1191.6316
$ws.Range("M74").Value = -317.6315999999999
$ws.Range("H77").Value = 2015.7333
$ws.Range("I77").Value = 1191.6316
$ws.Range("K77").Value = 5958.157999999999
$ws.Range("M77").Value = -1590.157999999999
$ws.Range("H97").Value = 1909.5
$ws.Range("I97").Value = 1108.4445
$ws.Range("K97").Value = 1108.4445
$ws.Range("M97").Value = -612.4445000000001
$ws.Range("H122").Value = 2279
$ws.Range("I122").Value = 2099
$ws.Range("K122").Value = 6297
$ws.Range("M122").Value = -3847
$ws.Range("H132").Value = 1611.2059
$ws.Range("I132").Value = 1436.5483
$ws.Range("J132").Value = 3416
$ws.Range("K132").Value = 4309.644899999999
$ws.Range("L132").Value = 10248
$ws.Range("M132").Value = -1779.644899999999
$ws.Range("N132").Value = -15308
$ws.Range("H136").Value = 6559.8
$ws.Range("I136").Value = 5800
$ws.Range("J136").Value = 7699.5
$ws.Range("K136").Value = 17400
$ws.Range("L136").Value = 23098.5
$ws.Range("M136").Value = -14850
$ws.Range("N136").Value = -28198.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2587.7
$ws.Range("I20").Value = 2541.889
$ws.Range("K20").Value = 2541.889
$ws.Range("M20").Value = -2294.889
$ws.Range("H94").Value = 17813.9
$ws.Range("I94").Value = 10180.8
$ws.Range("K94").Value = 10180.8
$ws.Range("M94").Value = -9729.799999999999
$ws.Range("H132").Value = 97260
$ws.Range("J132").Value = 97260
$ws.Range("L132").Value = 97260
$ws.Range("N132").Value = -107380
$ws.Range("H134").Value = 8922.647999999999
$ws.Range("I134").Value = 8463.034
$ws.Range("J134").Value = 10588.75
$ws.Range("K134").Value = 25389.102
$ws.Range("L134").Value = 31766.25
$ws.Range("M134").Value = -22854.102
$ws.Range("N134").Value = -36836.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3384.3215
$ws.Range("I31").Value = 3212.7693
$ws.Range("J31").Value = 3533
$ws.Range("K31").Value = 3212.7693
$ws.Range("L31").Value = 3533
$ws.Range("M31").Value = -2917.7693
$ws.Range("N31").Value = -4123
$ws.Range("H34").Value = 3384.3215
$ws.Range("I34").Value = 3212.7693
$ws.Range("J34").Value = 3533
$ws.Range("K34").Value = 3212.7693
$ws.Range("L34").Value = 3533
$ws.Range("M34").Value = -3010.7693
$ws.Range("N34").Value = -3937
$ws.Range("H132").Value = 2195.158
$ws.Range("I132").Value = 1837.9333
$ws.Range("K132").Value = 5513.7999
$ws.Range("M132").Value = -2983.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 473
$ws.Range("J24").Value = 567.5
$ws.Range("L24").Value = 1702.5
$ws.Range("N24").Value = -2162.5
$ws.Range("H109").Value = 4878.727
$ws.Range("J109").Value = 9735
$ws.Range("L109").Value = 29205
$ws.Range("N109").Value = -31285
$ws.Range("H140").Value = 10874883
$ws.Range("I140").Value = 13160700
$ws.Range("K140").Value = 39482100
$ws.Range("M140").Value = -39476920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 9998.333000000001
$ws.Range("I12").Value = 9998
$ws.Range("J12").Value = 9999
$ws.Range("K12").Value = 9998
$ws.Range("L12").Value = 9999
$ws.Range("M12").Value = -9858
$ws.Range("N12").Value = -10279
$ws.Range("H97").Value = 3290.3
$ws.Range("J97").Value = 7248.5
$ws.Range("L97").Value = 7248.5
$ws.Range("N97").Value = -8240.5
$ws.Range("H113").Value = 6216.091
$ws.Range("I113").Value = 4889
$ws.Range("K113").Value = 4889
$ws.Range("M113").Value = -2719

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H22").Value = 1214.1
$ws.Range("I22").Value = 902.9
$ws.Range("J22").Value = 1317.8334
$ws.Range("K22").Value = 902.9
$ws.Range("L22").Value = 1317.8334
$ws.Range("M22").Value = -607.9
$ws.Range("N22").Value = -1907.8334
$ws.Range("H27").Value = 1214.1
$ws.Range("I27").Value = 902.9
$ws.Range("J27").Value = 1317.8334
$ws.Range("K27").Value = 902.9
$ws.Range("L27").Value = 1317.8334
$ws.Range("M27").Value = -795.9
$ws.Range("N27").Value = -1531.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16976.445
$ws.Range("I122").Value = 4684.2856
$ws.Range("K122").Value = 14052.8568
$ws.Range("M122").Value = -11602.8568
